$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New IO-CE crosswalk rows to append (rows 260-300), restoring previously
# ignored sectors. Each entry is IO-sector (col A), CES category (col B),
# and optionally a second CES category (col C).
$rows = @(
  @("Religious organizations", "Fees and admissions"),
  @("Legal services", "Other household expenses"),
  @("Accounting, tax preparation, bookkeeping, and payroll services", "Other household expenses"),
  @("Specialized design services", "Other household expenses"),
  @("Scientific research and development services", "Other household expenses"),
  @("Advertising, public relations, and related services", "Other household expenses"),
  @("Funds, trusts, and other financial vehicles", "Other household expenses"),
  @("Securities and commodity contracts intermediation and brokerage", "Other household expenses"),
  @("Other financial investment activities", "Other household expenses"),
  @("Couriers and messengers", "Other household expenses"),
  @("Noncomparable imports", "Other household expenses"),
  @("Private households", "Other household expenses"),
  @("Other state and local government enterprises", "Other household expenses"),
  @("Individual and family services", "Other household expenses"),
  @("Other support services", "Other household expenses"),
  @("Veterinary services", "Other household expenses"),
  @("Employment services", "Other household expenses"),
  @("Business support services", "Other household expenses"),
  @("Travel arrangement and reservation services", "Other household expenses"),
  @("Investigation and security services", "Other household expenses"),
  @("Other real estate", "Mortgage interest and charges"),
  @("Pipeline transportation", "Fuel oil and other fuels"),
  @("Support activities for agriculture and forestry", "Water and other public services"),
  @("Grantmaking, giving, and social advocacy organizations", "Other household expenses"),
  @("Civic, social, professional, and similar organizations", "Other household expenses"),
  @("Commercial and industrial machinery and equipment rental and leasing", "Other household expenses"),
  @("Warehousing and storage", "Other household expenses"),
  @("Wholesale trade", "Other household expenses"),
  @("Other nonmetallic mineral mining and quarrying", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Iron and steel mills and ferroalloy manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Nonferrous metal (except copper and aluminum) rolling, drawing, extruding and alloying", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Nonferrous metal foundries", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Crown and closure manufacturing and metal stamping", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Plate work and fabricated structural product manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Metal can, box, and other metal container (light gauge) manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Hardware manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Spring and wire product manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Office machinery manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Metal cutting and forming machine tool manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Other engine equipment manufacturing", "Major appliances", "Small appliances, miscellaneous housewares"),
  @("Industrial gas manufacturing", "Major appliances", "Small appliances, miscellaneous housewares")
)

$startRow = 260
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $entry = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $entry[0]
  $ws.Cells.Item($r, 2).Value = $entry[1]
  if ($entry.Count -gt 2) {
    $ws.Cells.Item($r, 3).Value = $entry[2]
  }
}

$ws.Range("C300").Select() | Out-Null
